# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values on the 4d131514-... report row for the zh-cn and de-de sheets,
# reflecting a regenerated handback report with newer timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-11 12:28:37"
$wsZhCn.Range("H4").Value = "2016-03-11 12:28:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-11 12:28:40"
$wsDeDe.Range("H4").Value = "2016-03-11 12:29:00"
